$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: id_penyakit (A), penyakit (B), gejala (C) for 29 records (P0001..P0029)
$data = @(
    @('P0001', 'Rubella (campak)', 'Demam dan ruam yang menyebar ke seluruh tubuh'),
    @('P0002', 'Batuk rejan', 'Batuk keras, sesak napas'),
    @('P0003', 'Radang selaput otak', 'Sakit kepala, demam, leher kaku, flu'),
    @('P0004', 'Radang tenggorokan', 'Bersin, sulit menelan, air liur berlebihan, ruam, demam'),
    @('P0005', 'Demam berdarah', 'Ruam, kasar merah, radang tenggorokan, demam'),
    @('P0006', 'Sindrom Reye', 'Perilaku dramatis, kejang, koma'),
    @('P0007', 'Impetigo', 'Kulit melepuh, infeksi'),
    @('P0008', 'Kurap', 'Cincin merah pada kulit, kerontokan rambut kepala'),
    @('P0009', 'Lyme', 'Demam, menggigil, nyeri tubuh'),
    @('P0010', 'Flu', 'Demam tinggi, nyeri tubuh, kelelahan'),
    @('P0011', 'Alergi ', 'Bersin, mata berairi, hidung berair'),
    @('P0012', 'Radang amandel', 'Pembengkakan amandel'),
    @('P0013', 'Demam', 'Sakit tenggorokan, kedinginan, pipi merah'),
    @('P0014', 'Batuk', 'Flu, batuk'),
    @('P0015', 'Cacar air', 'Demam, sakit kepala, gatal, benjolan cairan muncul'),
    @('P0016', 'Panas dalam', 'Bibir pecah pecah, sakit, benjolan di bibir'),
    @('P0017', 'Infeksi saluran napas', 'Sesak, pneumonia, bronkitis, radang tenggorokan'),
    @('P0018', 'HIV/AIDS', 'Lemah, kekebalan tubuh berkurang'),
    @('P0019', 'Malaria', 'Demam, bintik merah, tulang nyeri'),
    @('P0020', 'Diare', 'Kolera, disentri, buang air besar sering, perut kembung'),
    @('P0021', 'TBC', 'Batuk, menggigil, demam, penurunan berat badan'),
    @('P0022', 'Campak', 'Diare, demam'),
    @('P0023', 'Tetanus', 'Infeksi luka'),
    @('P0024', 'Difteri', 'Demam, sakit tenggorokan, bercak putih di teggorokan'),
    @('P0025', 'Infeksi radang tenggorokan', 'Pilek, demam tinggi, bintik merah, nanah putih'),
    @('P0026', 'Rhinitis alergi', 'Hidung tersumbat, sakit tenggorokan, gatal, sakit kepala, nyeri wajah'),
    @('P0027', 'Infeksi telinga tengah', 'Demam, cairan bening mengalir dari telinga, sakit kepala'),
    @('P0028', 'Sesak napas', 'Wajah biru'),
    @('P0029', 'Kuning (jaundice)', 'Kulit kuning, demam')
)

$rowCount = $data.Count
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# Remove the now-obsolete last row (previously row 30)
$ws.Rows.Item($rowCount + 1).Delete() | Out-Null

# Resize column B to fit the longest "penyakit" entry
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# Select the header row, mirroring the author's final selection state
$ws.Range("A1:XFD1").Select() | Out-Null
